# cs-en-us-081pct.xlsx weekly refresh: roll the report forward one week
# (Vol 32 No 19 -> No 20; week of 5/5-5/11 -> 5/12-5/18) and drop in the
# newly collected crime-complaint figures for rows 15-30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: bump the issue number and the reporting week -----------
# A8 = "Volume 32   Number  19" -> only the trailing "19" run changes.
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "20"

# C9 = "Report Covering the Week  5/5/2025  Through  5/11/2025"
# Replace the *second* date first so the first replacement's length change
# (8 chars -> 9 chars) doesn't shift the second date's start offset.
$c9 = $ws.Range("C9")
$c9.Characters(46, 9).Text = "5/18/2025"
$c9.Characters(27, 8).Text = "5/12/2025"

# --- Column H got one tick narrower (bestFit shrank to match F/G) --------
$ws.Columns.Item(8).ColumnWidth = 5.43

# Row 15 (Rape) ------------------------------------------------------------
$ws.Range("N15").Value = -77.142857142857

# Row 16 (Robbery) ----------------------------------------------------------
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 10
$ws.Range("I16").Value = 45
$ws.Range("J16").Value = 46
$ws.Range("K16").Value = -2.173913043478
$ws.Range("L16").Value = -40.789473684210
$ws.Range("M16").Value = -56.310679611650
$ws.Range("N16").Value = -90.088105726872

# Row 17 (Fel. Assault) ------------------------------------------------------
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -37.5
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = -4.166666666666
$ws.Range("I17").Value = 90
$ws.Range("J17").Value = 115
$ws.Range("K17").Value = -21.739130434782
$ws.Range("L17").Value = -20.353982300885
$ws.Range("M17").Value = -15.094339622641
$ws.Range("N17").Value = -70.297029702970

# Row 18 (Burglary) -- D18/E18 flip from numbers to literal text ------------
$ws.Range("C18").Value = 3
$ws.Range("D14").Copy($ws.Range("D18"))        # D14 already holds text idx "0" @ style s=13
$ws.Range("E14").Copy($ws.Range("E18"))        # E14 already holds text idx "***.*" @ style s=13
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 100
$ws.Range("I18").Value = 29
$ws.Range("K18").Value = -21.621621621621
$ws.Range("L18").Value = -56.060606060606
$ws.Range("M18").Value = -59.154929577464
$ws.Range("N18").Value = -89.056603773584

# Row 19 (Gr. Larceny) -------------------------------------------------------
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 16.666666666666
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = -7.407407407407
$ws.Range("I19").Value = 89
$ws.Range("J19").Value = 105
$ws.Range("K19").Value = -15.238095238095
$ws.Range("L19").Value = -29.365079365079
$ws.Range("M19").Value = -10.101010101010
$ws.Range("N19").Value = -24.576271186440

# Row 20 (G.L.A.) -- C20 flips from text "0" to the number 1 -----------------
$ws.Range("I14").Copy($ws.Range("C20"))        # borrow numeric style s=14
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -60
$ws.Range("I20").Value = 31
$ws.Range("J20").Value = 33
$ws.Range("K20").Value = -6.060606060606
$ws.Range("L20").Value = -16.216216216216
$ws.Range("M20").Value = -13.888888888888
$ws.Range("N20").Value = -86.752136752136

# Row 21 (TOTAL, bold) -------------------------------------------------------
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = 11.764705882352
$ws.Range("F21").Value = 69
$ws.Range("G21").Value = 74
$ws.Range("H21").Value = -6.756756756756
$ws.Range("I21").Value = 293
$ws.Range("J21").Value = 341
$ws.Range("K21").Value = -14.076246334310
$ws.Range("L21").Value = -31.701631701631
$ws.Range("M21").Value = -32.488479262672
$ws.Range("N21").Value = -79.351656095842

# Row 22 (Transit) -- C22 flips from the number 1 to literal text "0" -------
$ws.Range("D22").Copy($ws.Range("C22"))        # D22 already holds text idx "0" @ style s=13

# Row 23 (Housing) -----------------------------------------------------------
$ws.Range("C23").Value = 4
$ws.Range("E23").Value = 300
$ws.Range("F23").Value = 12
$ws.Range("H23").Value = 33.333333333333
$ws.Range("I23").Value = 46
$ws.Range("J23").Value = 40
$ws.Range("K23").Value = 15
$ws.Range("L23").Value = 43.75
$ws.Range("M23").Value = 58.620689655172

# Row 24 (Petit Larceny) ------------------------------------------------------
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 70
$ws.Range("F24").Value = 54
$ws.Range("G24").Value = 47
$ws.Range("H24").Value = 14.893617021276
$ws.Range("I24").Value = 240
$ws.Range("J24").Value = 271
$ws.Range("K24").Value = -11.439114391143
$ws.Range("L24").Value = -24.528301886792
$ws.Range("M24").Value = -4.761904761904

# Row 25 (Retail Theft) -- C25 flips from the number 1 to literal text "0" --
$ws.Range("D25").Value = 1
$ws.Range("G22").Copy($ws.Range("C25"))        # G22 already holds text idx "0" @ style s=13
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 5
$ws.Range("G25").Value = 8
$ws.Range("H25").Value = -37.5
$ws.Range("I25").Value = 23
$ws.Range("J25").Value = 43
$ws.Range("K25").Value = -46.511627906976
$ws.Range("L25").Value = -78.301886792452

# Row 26 (Misd. Assault) ------------------------------------------------------
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = 28.571428571428
$ws.Range("F26").Value = 34
$ws.Range("G26").Value = 42
$ws.Range("H26").Value = -19.047619047619
$ws.Range("I26").Value = 145
$ws.Range("J26").Value = 152
$ws.Range("K26").Value = -4.605263157894
$ws.Range("L26").Value = -17.142857142857
$ws.Range("M26").Value = -44.866920152091

# Row 27 (UCR Rape*) -- G27/H27 flip from numbers to literal text -----------
$ws.Range("C27").Copy($ws.Range("G27"))        # C27 already holds text idx "0" @ style s=13
$ws.Range("E27").Copy($ws.Range("H27"))        # E27 already holds text idx "***.*" @ style s=13

# Row 28 (Other Sex Crimes) ---------------------------------------------------
$ws.Range("D28").Value = 3
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = -66.666666666666
$ws.Range("J28").Value = 19
$ws.Range("K28").Value = -15.789473684210

# Row 29 (Shooting Vic.) -------------------------------------------------------
$ws.Range("F29").Value = 1
$ws.Range("M29").Value = -65
$ws.Range("N29").Value = -91.463414634146

# Row 30 (Shooting Inc.) -------------------------------------------------------
$ws.Range("F30").Value = 1
$ws.Range("M30").Value = -66.666666666666
$ws.Range("N30").Value = -93.150684931506
